$d = $word.ActiveDocument

$newText = "Dates de la campanya Pegasus: Del 8 al 17 d" + [char]0x2019 + "octubre, del 7 al 16 de novembre,"

$oldText1 = " Dates de la campanya 2018 en què usem la constel·lació Perseus 30 d'octubre al novembre 8 i 29 de novembre de desembre 8"
$oldText2 = "Dates de la campanya 2018 en què usem la constel·lació Perseus 30 d'octubre al novembre 8 i 29 de novembre de desembre 8"

# Occurrence 1 (paragraph with leading red space run)
$rng = $d.Content
$rng.Find.Execute($oldText1, $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
$rng.Delete()
$rng.InsertAfter($newText)

# Occurrences 2-4 (same simpler pattern, repeated 3 times)
for ($i = 0; $i -lt 3; $i++) {
    $rng2 = $d.Content
    $rng2.Find.Execute($oldText2, $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
    $rng2.Delete()
    $rng2.InsertAfter($newText)
}
